$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 34.76923
$ws.Range("I11").Value = 34.76923
$ws.Range("K11").Value = 34.76923
$ws.Range("M11").Value = 105.23077
$ws.Range("H33").Value = 6399.0625
$ws.Range("I33").Value = 88.5
$ws.Range("K33").Value = 88.5
$ws.Range("M33").Value = 140.5
$ws.Range("H64").Value = 3216.25
$ws.Range("I64").Value = 3076.5715
$ws.Range("J64").Value = 3542.1667
$ws.Range("K64").Value = 3076.5715
$ws.Range("L64").Value = 3542.1667
$ws.Range("M64").Value = -2828.5715
$ws.Range("N64").Value = -4038.1667
$ws.Range("H67").Value = 3216.25
$ws.Range("I67").Value = 3076.5715
$ws.Range("J67").Value = 3542.1667
$ws.Range("K67").Value = 3076.5715
$ws.Range("L67").Value = 3542.1667
$ws.Range("M67").Value = -2218.5715
$ws.Range("N67").Value = -5258.1667
$ws.Range("H115").Value = 709.6316
$ws.Range("I115").Value = 452.91666
$ws.Range("J115").Value = 1149.7142
$ws.Range("K115").Value = 1358.74998
$ws.Range("L115").Value = 3449.1426
$ws.Range("M115").Value = 208.2500199999999
$ws.Range("N115").Value = -6583.142599999999
$ws.Range("H137").Value = 1090.6428
$ws.Range("I137").Value = 928.4737
$ws.Range("J137").Value = 1433
$ws.Range("K137").Value = 2785.4211
$ws.Range("L137").Value = 4299
$ws.Range("M137").Value = -235.4211
$ws.Range("N137").Value = -9399
$ws.Range("H138").Value = 2582.875
$ws.Range("I138").Value = 1809.825
$ws.Range("J138").Value = 4515.5
$ws.Range("K138").Value = 5429.475
$ws.Range("L138").Value = 13546.5
$ws.Range("M138").Value = -289.4750000000004
$ws.Range("N138").Value = -23826.5
$ws.Range("H141").Value = 1835.6
$ws.Range("I141").Value = 1242.8125
$ws.Range("J141").Value = 2889.4443
$ws.Range("K141").Value = 3728.4375
$ws.Range("L141").Value = 8668.332900000001
$ws.Range("M141").Value = 1451.5625
$ws.Range("N141").Value = -19028.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3190.923
$ws.Range("I61").Value = 3210.353
$ws.Range("J61").Value = 3058.8
$ws.Range("K61").Value = 3210.353
$ws.Range("L61").Value = 3058.8
$ws.Range("M61").Value = -2998.353
$ws.Range("N61").Value = -3482.8
$ws.Range("H74").Value = 1137.8182
$ws.Range("I74").Value = 627
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 627
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = 247
$ws.Range("N74").Value = -4248
$ws.Range("H76").Value = 29666.666
$ws.Range("J76").Value = 29666.666
$ws.Range("L76").Value = 29666.666
$ws.Range("N76").Value = -30342.666
$ws.Range("H77").Value = 1137.8182
$ws.Range("I77").Value = 627
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 3135
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = 1233
$ws.Range("N77").Value = -21236
$ws.Range("H79").Value = 29666.666
$ws.Range("J79").Value = 29666.666
$ws.Range("L79").Value = 29666.666
$ws.Range("N79").Value = -32006.666
$ws.Range("H132").Value = 2572.5789
$ws.Range("I132").Value = 2377.926
$ws.Range("J132").Value = 3050.3635
$ws.Range("K132").Value = 7133.778
$ws.Range("L132").Value = 9151.0905
$ws.Range("M132").Value = -4603.778
$ws.Range("N132").Value = -14211.0905
$ws.Range("H136").Value = 3190.923
$ws.Range("I136").Value = 3210.353
$ws.Range("J136").Value = 3058.8
$ws.Range("K136").Value = 9631.059000000001
$ws.Range("L136").Value = 9176.400000000001
$ws.Range("M136").Value = -7081.059000000001
$ws.Range("N136").Value = -14276.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5278.2812
$ws.Range("I31").Value = 4214.72
$ws.Range("J31").Value = 5960.0513
$ws.Range("K31").Value = 4214.72
$ws.Range("L31").Value = 5960.0513
$ws.Range("M31").Value = -3919.72
$ws.Range("N31").Value = -6550.0513
$ws.Range("H34").Value = 5278.2812
$ws.Range("I34").Value = 4214.72
$ws.Range("J34").Value = 5960.0513
$ws.Range("K34").Value = 4214.72
$ws.Range("L34").Value = 5960.0513
$ws.Range("M34").Value = -4012.72
$ws.Range("N34").Value = -6364.0513
$ws.Range("H51").Value = 9599
$ws.Range("J51").Value = 9599
$ws.Range("L51").Value = 9599
$ws.Range("N51").Value = -11071
$ws.Range("H60").Value = 8251
$ws.Range("J60").Value = 8251
$ws.Range("L60").Value = 8251
$ws.Range("N60").Value = -9273
$ws.Range("H61").Value = 9599
$ws.Range("J61").Value = 9599
$ws.Range("L61").Value = 9599
$ws.Range("N61").Value = -10295
$ws.Range("H99").Value = 1966.0646
$ws.Range("I99").Value = 1824.3914
$ws.Range("J99").Value = 2373.375
$ws.Range("K99").Value = 1824.3914
$ws.Range("L99").Value = 2373.375
$ws.Range("M99").Value = -326.3914
$ws.Range("N99").Value = -5369.375
$ws.Range("H126").Value = 1966.0646
$ws.Range("I126").Value = 1824.3914
$ws.Range("J126").Value = 2373.375
$ws.Range("K126").Value = 5473.174199999999
$ws.Range("L126").Value = 7120.125
$ws.Range("M126").Value = -3003.174199999999
$ws.Range("N126").Value = -12060.125
$ws.Range("H132").Value = 1829.0264
$ws.Range("I132").Value = 1480.5555
$ws.Range("J132").Value = 2684.3635
$ws.Range("K132").Value = 4441.666499999999
$ws.Range("L132").Value = 8053.0905
$ws.Range("M132").Value = -1911.666499999999
$ws.Range("N132").Value = -13113.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1118.091
$ws.Range("I5").Value = 630.7273
$ws.Range("J5").Value = 1605.4546
$ws.Range("K5").Value = 1892.1819
$ws.Range("L5").Value = 4816.3638
$ws.Range("M5").Value = -1780.1819
$ws.Range("N5").Value = -5040.3638
$ws.Range("H56").Value = 3944406.2
$ws.Range("I56").Value = 3944406.2
$ws.Range("K56").Value = 3944406.2
$ws.Range("M56").Value = -3943876.2
$ws.Range("H93").Value = 3980
$ws.Range("J93").Value = 3980
$ws.Range("L93").Value = 11940
$ws.Range("N93").Value = -15684
$ws.Range("H113").Value = 1016.45056
$ws.Range("I113").Value = 925.75
$ws.Range("J113").Value = 1020.62067
$ws.Range("K113").Value = 2777.25
$ws.Range("L113").Value = 3061.86201
$ws.Range("M113").Value = -607.25
$ws.Range("N113").Value = -7401.86201
$ws.Range("H131").Value = 759.5625
$ws.Range("I131").Value = 430
$ws.Range("J131").Value = 806.6429000000001
$ws.Range("K131").Value = 1290
$ws.Range("L131").Value = 2419.9287
$ws.Range("M131").Value = 3750
$ws.Range("N131").Value = -12499.9287
$ws.Range("H135").Value = 1118.091
$ws.Range("I135").Value = 630.7273
$ws.Range("J135").Value = 1605.4546
$ws.Range("K135").Value = 5676.545700000001
$ws.Range("L135").Value = 14449.0914
$ws.Range("M135").Value = -3141.545700000001
$ws.Range("N135").Value = -19519.0914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 372296.28
$ws.Range("H30").Value = 372296.28
$ws.Range("H132").Value = 2572.1714
$ws.Range("I132").Value = 1915.3914
$ws.Range("J132").Value = 3831
$ws.Range("K132").Value = 5746.174199999999
$ws.Range("L132").Value = 11493
$ws.Range("M132").Value = -3216.174199999999
$ws.Range("N132").Value = -16553

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 718.8788
$ws.Range("I22").Value = 543.1818
$ws.Range("J22").Value = 1070.2727
$ws.Range("K22").Value = 543.1818
$ws.Range("L22").Value = 1070.2727
$ws.Range("M22").Value = -248.1818
$ws.Range("N22").Value = -1660.2727
$ws.Range("H27").Value = 718.8788
$ws.Range("I27").Value = 543.1818
$ws.Range("J27").Value = 1070.2727
$ws.Range("K27").Value = 543.1818
$ws.Range("L27").Value = 1070.2727
$ws.Range("M27").Value = -436.1818
$ws.Range("N27").Value = -1284.2727
$ws.Range("H68").Value = 2005.6171
$ws.Range("I68").Value = 1993.3334
$ws.Range("J68").Value = 2034.5714
$ws.Range("K68").Value = 1993.3334
$ws.Range("L68").Value = 2034.5714
$ws.Range("M68").Value = -1244.3334
$ws.Range("N68").Value = -3532.5714
$ws.Range("H71").Value = 2005.6171
$ws.Range("I71").Value = 1993.3334
$ws.Range("J71").Value = 2034.5714
$ws.Range("K71").Value = 9966.666999999999
$ws.Range("L71").Value = 10172.857
$ws.Range("M71").Value = -6222.666999999999
$ws.Range("N71").Value = -17660.857
$ws.Range("H100").Value = 3355.818
$ws.Range("I100").Value = 3164.8333
$ws.Range("J100").Value = 3585
$ws.Range("K100").Value = 3164.8333
$ws.Range("L100").Value = 3585
$ws.Range("M100").Value = -2623.8333
$ws.Range("N100").Value = -4667
$ws.Range("H132").Value = 7814.364
$ws.Range("I132").Value = 9770.286
$ws.Range("J132").Value = 4391.5
$ws.Range("K132").Value = 29310.858
$ws.Range("L132").Value = 13174.5
$ws.Range("M132").Value = -26780.858
$ws.Range("N132").Value = -18234.5
